$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 1.28
$ws.Range("F3").Value = 1.19
$ws.Range("E4").Value = 1.25
$ws.Range("C5").Value = 1.36
$ws.Range("C6").Value = 1.5
$ws.Range("G6").Value = 1
$ws.Range("E7").Value = 1.93
